# Apply cryptos list price/volume update (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.937.62"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "3.542.82"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "618.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.83"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "3.537.87"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.20"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "4.113.65"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.46"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "612.84"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "70.947.12"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.530.33"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.84"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.84%  "
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.04"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.72%  "
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.31"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.87"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.16"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.04"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.19"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.32%  "
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.87"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "605.64"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.91%  "
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.52"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0477"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.99"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("E41").Value = "  +3.25%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "3.382.33"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "32.33"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.59"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.16%  "
